$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 133; this pushes the existing rows 133-242
# down to 134-243 (matching the shift seen throughout the diff) and grows
# the sheet's used range from A1:R242 to A1:R243.
$ws.Rows.Item(133).EntireRow.Insert()

# Populate the newly inserted row 133 with the new price-report record.
$ws.Cells.Item(133, 1).Value = 4
$ws.Cells.Item(133, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(133, 3).Value = "Los Lagos"
$ws.Cells.Item(133, 4).Value = 44651
$ws.Cells.Item(133, 5).Value = 10
$ws.Cells.Item(133, 6).Value = 100112003
$ws.Cells.Item(133, 7).Value = "Ajo"
$ws.Cells.Item(133, 8).Value = "Chino"
$ws.Cells.Item(133, 9).Value = "Primera"
$ws.Cells.Item(133, 10).Value = 35
$ws.Cells.Item(133, 11).Value = 21000
$ws.Cells.Item(133, 12).Value = 21000
$ws.Cells.Item(133, 13).Value = 21000
$ws.Cells.Item(133, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(133, 15).Value = "China"
$ws.Cells.Item(133, 16).Value = 2100
$ws.Cells.Item(133, 17).Value = 10
$ws.Cells.Item(133, 18).Value = "Hortaliza"
